# Revise responsive design implementation — append the next daily log
# row (row 68) to each of the four data sheets, matching the style /
# number-format already used by the existing rows (date in col A,
# hex-string columns B-E, decimal columns F-I).

$wb = $excel.ActiveWorkbook

$dateFormat = "YYYY-MM-DD HH:MM:SS"
$newDate = [double]"45854.43328703703"

$rows = @{
    "DE_LFT_#1" = @{
        A = $newDate
        B = "0x01,0x7c"
        C = "0x00,0xa6,0x46,0x93,0x3c,0x23,0x3f,0x43,0xe8,0xa0,"
        D = "0x01,0x4C"
        E = "0x14"
        F = 380
        G = [double]"7.598631275147109e+23"
        H = 332
        I = 14
    }
    "DE_LFT_#2" = @{
        A = $newDate
        B = "0x01,0x7c"
        C = "0x00,0xa6,0x60,0x33,0x96,0x39,0x62,0xd0,0x5e,0x78,"
        D = "0x01,0x50"
        E = "0xe"
        F = 380
        G = [double]"5.68432987514711e+23"
        H = 336
        I = 14
    }
    "DE_PLT_#1" = @{
        A = $newDate
        B = "0x00,0x82"
        C = "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x0b,0x40,0x0c,"
        D = "0x00,0x7B"
        E = "0x7"
        F = 130
        G = [double]"5.68631262647114e+23"
        H = 123
        I = 7
    }
    "DE_PLT_#2" = @{
        A = $newDate
        B = "0x00,0x82"
        C = "0xd0,0x97,0x78,0x01,0x00,0x00,0x0e,0x3f,0x0c,0x0c,"
        D = "0x00,0x7B"
        E = "0x3"
        F = 130
        G = [double]"9.85046333984776e+23"
        H = 123
        I = 3
    }
}

foreach ($sheetName in $rows.Keys) {
    $ws = $wb.Worksheets.Item($sheetName)
    $data = $rows[$sheetName]

    $ws.Range("A68").Value = $data.A
    $ws.Range("A68").NumberFormat = $dateFormat

    $ws.Range("B68").Value = $data.B
    $ws.Range("C68").Value = $data.C
    $ws.Range("D68").Value = $data.D
    $ws.Range("E68").Value = $data.E
    $ws.Range("F68").Value = $data.F
    $ws.Range("G68").Value = $data.G
    $ws.Range("H68").Value = $data.H
    $ws.Range("I68").Value = $data.I
}
